$wb = $excel.ActiveWorkbook

# ALC row 5: Met a Sticky End | Animal Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 289
$ws.Range("I5").Value = 204.6
$ws.Range("K5").Value = 204.6
$ws.Range("M5").Value = -89.59999999999999

# ALC row 58: A Matter of Vital Importance | Mega-Potion of Vitality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1307.8667
$ws.Range("J58").Value = 2261.875
$ws.Range("L58").Value = 6785.625
$ws.Range("N58").Value = -7085.625

# ALC row 70: Consecrating Congregation | Holy Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 3000
$ws.Range("K70").Value = 9000
$ws.Range("M70").Value = -8730

# ALC row 73: Curbing the Contagion (L) | Holy Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 3000
$ws.Range("K73").Value = 9000
$ws.Range("M73").Value = -8064

# ALC row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4123.2666
$ws.Range("I86").Value = 3209.9
$ws.Range("J86").Value = 5950
$ws.Range("K86").Value = 3209.9
$ws.Range("L86").Value = 5950
$ws.Range("M86").Value = -2086.9
$ws.Range("N86").Value = -8196

# ALC row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4123.2666
$ws.Range("I89").Value = 3209.9
$ws.Range("J89").Value = 5950
$ws.Range("K89").Value = 16049.5
$ws.Range("L89").Value = 29750
$ws.Range("M89").Value = -10433.5
$ws.Range("N89").Value = -40982

# ALC row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7938590
$ws.Range("I132").Value = 8335344.5
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 25006033.5
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -25003503.5
$ws.Range("N132").Value = -15560

# ARM row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1014.5217
$ws.Range("I61").Value = 821.05
$ws.Range("K61").Value = 821.05
$ws.Range("M61").Value = -609.05

# ARM row 63: Rivets Run through It | Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2189.9
$ws.Range("I63").Value = 2237.375
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2237.375
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1551.375
$ws.Range("N63").Value = -3372

# ARM row 66: A Riveting Revival (L) | Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2189.9
$ws.Range("I66").Value = 2237.375
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 11186.875
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -7754.875
$ws.Range("N66").Value = -16864

# ARM row 96: The Gauntlet Is Cast | High Steel Gauntlets of Fending
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 18335.75
$ws.Range("J96").Value = 18335.75
$ws.Range("L96").Value = 18335.75
$ws.Range("N96").Value = -23827.75

# ARM row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1014.5217
$ws.Range("I136").Value = 821.05
$ws.Range("K136").Value = 2463.15
$ws.Range("M136").Value = 86.85000000000036

# BSM row 35: Lancers' Creed | Crowsbeak Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("N35").Value = 0

# BSM row 82: Spirituality Inspector | Titanium Lump Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 29063.75
$ws.Range("J82").Value = 36000
$ws.Range("L82").Value = 36000
$ws.Range("N82").Value = -36766

# BSM row 85: The Clamor for Hammers (L) | Titanium Lump Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 29063.75
$ws.Range("J85").Value = 36000
$ws.Range("L85").Value = 36000
$ws.Range("N85").Value = -38652

# BSM row 86: Through Thick and Thin | Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2996.5527
$ws.Range("I86").Value = 3311.7036
$ws.Range("J86").Value = 2223
$ws.Range("K86").Value = 3311.7036
$ws.Range("L86").Value = 2223
$ws.Range("M86").Value = -2188.7036
$ws.Range("N86").Value = -4469

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2996.5527
$ws.Range("I89").Value = 3311.7036
$ws.Range("J89").Value = 2223
$ws.Range("K89").Value = 16558.518
$ws.Range("L89").Value = 11115
$ws.Range("M89").Value = -10942.518
$ws.Range("N89").Value = -22347

# CRP row 97: Wood That You Could | Larch Bracelets
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("N97").Value = 0

# CUL row 5: What a Sap | Maple Syrup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1522.25
$ws.Range("I5").Value = 1522.25
$ws.Range("K5").Value = 4566.75
$ws.Range("M5").Value = -4454.75

# CUL row 131: The Mountain Steeped | Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 897.79
$ws.Range("J131").Value = 935.043
$ws.Range("L131").Value = 2805.129
$ws.Range("N131").Value = -12885.129

# CUL row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1522.25
$ws.Range("I135").Value = 1522.25
$ws.Range("K135").Value = 13700.25
$ws.Range("M135").Value = -11165.25

# CUL row 136: Simple Is Hardest | Spaghetti al Olio e Peperoncino
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1799.9166
$ws.Range("I136").Value = 1066.6666
$ws.Range("J136").Value = 3999.6667
$ws.Range("K136").Value = 3199.9998
$ws.Range("L136").Value = 11999.0001
$ws.Range("M136").Value = 1900.0002
$ws.Range("N136").Value = -22199.0001

# GSM row 80: Needs More Prayerbell | Hardsilver Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4887.143
$ws.Range("I80").Value = 3652
$ws.Range("J80").Value = 5381.2
$ws.Range("K80").Value = 3652
$ws.Range("L80").Value = 5381.2
$ws.Range("M80").Value = -2654
$ws.Range("N80").Value = -7377.2

# GSM row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4887.143
$ws.Range("I83").Value = 3652
$ws.Range("J83").Value = 5381.2
$ws.Range("K83").Value = 18260
$ws.Range("L83").Value = 26906
$ws.Range("M83").Value = -13268
$ws.Range("N83").Value = -36890

# LTW row 40: Best Served Toad | Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2983.8
$ws.Range("I40").Value = 2883.8333
$ws.Range("J40").Value = 3133.75
$ws.Range("K40").Value = 2883.8333
$ws.Range("L40").Value = 3133.75
$ws.Range("M40").Value = -2747.8333
$ws.Range("N40").Value = -3405.75

# LTW row 100: Tiger in the Sack | Tiger Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3300
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

# WVR row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 700
$ws.Range("I81").Value = 700
$ws.Range("K81").Value = 1400
$ws.Range("M81").Value = -339

# WVR row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 700
$ws.Range("I84").Value = 700
$ws.Range("K84").Value = 7000
$ws.Range("M84").Value = -1696

# WVR row 108: Lovely Leggings | Brightlinen Bottoms of Striking
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0

# WVR row 110: Suits You | Iridescent Acton of Aiming
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 17822
$ws.Range("J110").Value = 17822
$ws.Range("L110").Value = 17822
$ws.Range("N110").Value = -26002

# WVR row 125: Color Coated | Almasty Serge Coat of Healing
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 78326.664
$ws.Range("J125").Value = 78326.664
$ws.Range("L125").Value = 78326.664
$ws.Range("N125").Value = -88166.664

# WVR row 126: A Polished Purchase | Snow Linen
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 100000960
$ws.Range("I126").Value = 100000960
$ws.Range("K126").Value = 300002880
$ws.Range("M126").Value = -300000410

# WVR row 131: A Better Bottom Line | AR-Caean Velvet Bottoms of Scouting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 79995
$ws.Range("J131").Value = 79995
$ws.Range("L131").Value = 79995
$ws.Range("N131").Value = -90075

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3071.2917
$ws.Range("I132").Value = 2632.1875
$ws.Range("K132").Value = 7896.5625
$ws.Range("M132").Value = -5366.5625
